$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from an existing "Manual Status = Passed" cell (G3)
# so the same cell style is reused rather than a brand-new one created.
$ws.Range("G3").Copy()

$rows = @(25, 37, 38, 39, 40, 41)
foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $cell.PasteSpecial(-4122)  # xlPasteFormats
    $cell.Value = "Passed"
}

$excel.CutCopyMode = $false

# Match the scrolled viewport / selection left behind by the edit.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("G41").Select()
